# Apply the UserSkills workbook edits:
#  - POST sheet: scenario id US15 -> US17
#  - DELETE sheet: status message typo fix "User Skill Map Not Found"
#                  -> "User skill Map Not Found" (rows 3-6)
#  - PUT sheet becomes the active/selected sheet (was DELETE before)

$wb = $excel.ActiveWorkbook

$wsPost   = $wb.Worksheets.Item(2)
$wsPut    = $wb.Worksheets.Item(3)
$wsDelete = $wb.Worksheets.Item(4)

# --- POST sheet: update the UserSkills_ID value in row 2 ---
$wsPost.Range("B2").Value = "US17"

# --- DELETE sheet: fix the StatusMessage text in rows 3-6 ---
$wsDelete.Range("D3").Value = "User skill Map Not Found"
$wsDelete.Range("D4").Value = "User skill Map Not Found"
$wsDelete.Range("D5").Value = "User skill Map Not Found"
$wsDelete.Range("D6").Value = "User skill Map Not Found"

# --- Switch the active tab from DELETE to PUT, update selection ---
$wsPut.Activate() | Out-Null
$wsPut.Range("E3").Select() | Out-Null
$wsDelete.Range("E6").Select() | Out-Null

$wsPut.Activate() | Out-Null
